$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (rows 2-10) "Förändrad" date value from 45221 to 45224
for ($r = 2; $r -le 10; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45221) {
        $cell.Value2 = 45224
    }
}
